$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) from the last existing row down to the new rows
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 8
$ws.Range("A8").Value = 4021
$ws.Range("B8").Value = "presunto oval"
$ws.Range("C8").Value = 412

# Row 9
$ws.Range("A9").Value = 5023
$ws.Range("B9").Value = "bacon fracionado"
$ws.Range("C9").Value = 87

# Row 10
$ws.Range("A10").Value = 5009
$ws.Range("B10").Value = "bacon manta"
$ws.Range("C10").Value = 200

# Row 11
$ws.Range("A11").Value = 4016
$ws.Range("B11").Value = "cuzcuz de flocao"
$ws.Range("C11").Value = 32

# Row 12
$ws.Range("A12").Value = 8005
$ws.Range("B12").Value = "teste "
$ws.Range("C12").Value = 120

# Update the active selection to reflect the new last cell, as in the final workbook
$ws.Range("C13").Select()
